$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# B11 currently holds the shared string "Good Night" (text). The target
# change replaces its value with the text "1" — a brand-new shared string,
# stored as text (not a number), while leaving the cell's existing style
# untouched.
#
# Simply assigning Range.Value = "1" would make Excel infer a numeric
# value (losing the "t=s" shared-string/text typing), and forcing text via
# NumberFormat "@" or a leading apostrophe (quote-prefix) both end up
# mutating/creating a cell style, which the target workbook does not do.
#
# Instead: build the text value as a formula result (="1"), which yields a
# genuine text cell, then copy/paste only the *value* into B11 so its
# existing formatting/style is left completely untouched. A scratch cell
# far outside the worksheet's used range is used as scratch space and
# cleared afterwards so nothing else is left behind.

$scratch = $ws.Range("Z100")
$scratch.Formula = "=""1"""
$scratch.Copy()
$ws.Range("B11").PasteSpecial(-4163)
$scratch.Clear()
